# Update Schwall_Gayathri_Timesheet_october_2021 - Copy - Copy.xlsx
# Fill in timesheet rows 25-30 on the "Week1" sheet with new task entries,
# and move the sheet's viewport/selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week1")

# Row 25: SlNo 9, Date 13-Oct-2021, Task "camunda Training", Desc "camunda,workflow", Hrs 2
$ws.Range("B25").Value = 9
$ws.Range("C25").Value = 44482
$ws.Range("D25").Value = "camunda Training"
$ws.Range("E25").Value = "camunda,workflow"
$ws.Range("F25").Value = 2

# Row 26: Task "camunda", Desc "camunda springboot installation", Hrs 1
$ws.Range("D26").Value = "camunda"
$ws.Range("E26").Value = "camunda springboot installation"
$ws.Range("F26").Value = 1

# Row 27: Task "java", Desc "multithreading", Hrs 1
$ws.Range("D27").Value = "java"
$ws.Range("E27").Value = "multithreading"
$ws.Range("F27").Value = 1

# Row 28: Task "java Task", Desc "Multithreading", Hrs 2
$ws.Range("D28").Value = "java Task"
$ws.Range("E28").Value = "Multithreading"
$ws.Range("F28").Value = 2

# Row 29: Task "camunda", Desc "spring boot setup,Tomcat setup", Hrs 1
$ws.Range("D29").Value = "camunda"
$ws.Range("E29").Value = "spring boot setup,Tomcat setup"
$ws.Range("F29").Value = 1

# Row 30: Task "revision", Desc "overall revision camunda", Hrs 3
$ws.Range("D30").Value = "revision"
$ws.Range("E30").Value = "overall revision camunda"
$ws.Range("F30").Value = 3

# Update view: scroll so row 4 is at top, and select F30 as the active cell.
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("F30").Select()
